# Update position of Flash
# - Sheet2: add "Value" (col E) and "Package" (col F) BOM columns for each
#   placement row (rows 1-32), sourced from the component datasheet info.
# - Sheet1: move the active-cell selection.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# Designator -> (Value, Package) for each row, in top-to-bottom sheet order.
# Writing E before F, row by row, keeps the shared-string table insertion
# order identical to the source edit.
$rows = @(
    @{ Row = 1;  E = "0.1uF";       F = 1206 }              # C1
    @{ Row = 2;  E = "22uF";        F = 1206 }              # C10
    @{ Row = 3;  E = "0.1uF";       F = 1206 }              # C11
    @{ Row = 4;  E = "0.1uF";       F = 1206 }              # C12
    @{ Row = 5;  E = "10uF";        F = 1206 }              # C13
    @{ Row = 6;  E = "10uF";        F = 1206 }              # C14
    @{ Row = 7;  E = "0.1uF";       F = 1206 }              # C2
    @{ Row = 8;  E = "0.1uF";       F = 1206 }              # C3
    @{ Row = 9;  E = "0.1uF";       F = 1206 }              # C4
    @{ Row = 10; E = "100uF";       F = 1206 }              # C5
    @{ Row = 11; E = "22uF";        F = 1206 }              # C6
    @{ Row = 12; E = "100uF";       F = 1206 }              # C9
    @{ Row = 13; E = $null;         F = "LED-1206" }        # D1
    @{ Row = 14; E = $null;         F = "LED-1206" }        # D2
    @{ Row = 15; E = $null;         F = "LED-1206" }        # D3
    @{ Row = 16; E = $null;         F = "LED-1206" }        # D4
    @{ Row = 17; E = "TXB0108PWR";  F = "SOP65P640X120-20N" } # IC1
    @{ Row = 18; E = "TXB0108PWR";  F = "SOP65P640X120-20N" } # IC2
    @{ Row = 19; E = "AMS1117-1.8"; F = "SOT229P700X180-4N" } # IC3
    @{ Row = 20; E = "AMS1117-1.8"; F = "SOT229P700X180-4N" } # IC5
    @{ Row = 21; E = "2N7002K-7";   F = "SOT96P240X100-3N" }  # Q5
    @{ Row = 22; E = "2N7002K-7";   F = "SOT96P240X100-3N" }  # Q6
    @{ Row = 23; E = "2N7002K-7";   F = "SOT96P240X100-3N" }  # Q7
    @{ Row = 24; E = "2N7002K-7";   F = "SOT96P240X100-3N" }  # Q8
    @{ Row = 25; E = "10k";         F = 1206 }              # R10
    @{ Row = 26; E = 330;           F = 1206 }              # R11
    @{ Row = 27; E = "10k";         F = 1206 }              # R12
    @{ Row = 28; E = 330;           F = 1206 }              # R13
    @{ Row = 29; E = "10k";         F = 1206 }              # R14
    @{ Row = 30; E = 330;           F = 1206 }              # R15
    @{ Row = 31; E = "10k";         F = 1206 }              # R16
    @{ Row = 32; E = 330;           F = 1206 }              # R17
)

foreach ($r in $rows) {
    if ($null -ne $r.E) {
        $ws2.Range("E$($r.Row)").Value = $r.E
    }
    if ($null -ne $r.F) {
        $ws2.Range("F$($r.Row)").Value = $r.F
    }
}

# Sheet1: move selection from J9 to H8.
$ws1 = $wb.Worksheets.Item("Sheet1")
$null = $ws1.Range("H8").Select()
